# "added printable to api"
# The inspection data for the 2nd sample record (11 Hay Ln / Mike Hamilton /
# Carlos Amana) is promoted into rows 3 & 4 (replacing the old sample rows
# 3-6), dates/labels are tidied up (trailing spaces / fuller date), and short
# free-text notes are now written into the corresponding Mold/Electrical/
# Foundation note columns (N, P, S) for the "printable" report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two extra sample rows (5 & 6) - only rows 1, 3 and 4 remain afterwards.
$ws.Rows.Item(5).EntireRow.Delete()
$ws.Rows.Item(5).EntireRow.Delete()

# --- Row 3 -------------------------------------------------------------
$ws.Range("A3").Value = "11 Hay Ln "
$ws.Range("B3").Value = "12/18/2022 "
$ws.Range("C3").Value = "Mike Hamilton "
$ws.Range("D3").Value = "Carlos Amana"
$ws.Range("E3").Value = "1087 Express Drive N "
$ws.Range("F3").Value = "555-294-3398 "
$ws.Range("G3").Value = "Yes"
$ws.Range("H3").Value = "No"
$ws.Range("I3").Value = "Yes"
$ws.Range("J3").Value = "No"
$ws.Range("K3").Value = "No"
$ws.Range("L3").Value = "Yes"
$ws.Range("M3").Value = "No"
$ws.Range("N3").Value = "Primary bathroom has black mold behind the sink. Requires complete remodel "
$ws.Range("P3").Value = "Wiring is old and ungrounded, should rewire entire home "
$ws.Range("S3").Value = "Big crack in Foundation "

# --- Row 4 (mirrors row 3, but without a Framing/Roof note) -------------
$ws.Range("A4").Value = "11 Hay Ln "
$ws.Range("B4").Value = "12/18/2022 "
$ws.Range("C4").Value = "Mike Hamilton "
$ws.Range("D4").Value = "Carlos Amana"
$ws.Range("E4").Value = "1087 Express Drive N "
$ws.Range("F4").Value = "555-294-3398 "
$ws.Range("G4").Value = "Yes"
$ws.Range("H4").Value = "No"
$ws.Range("I4").Value = "Yes"
$ws.Range("J4").Value = "No"
$ws.Range("K4").Value = "No"
$ws.Range("L4").Value = "Yes"
$ws.Range("M4").Value = "No"
$ws.Range("N4").Value = "Primary bathroom has black mold behind the sink. Requires complete remodel "
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = "Wiring is old and ungrounded, should rewire entire home "
$ws.Range("R4").ClearContents()
$ws.Range("S4").Value = "Big crack in Foundation "

# Columns whose content footprint changed get re-sized to fit the new
# (shorter) data set; columns that weren't touched by the edit keep Excel's
# original best-fit widths automatically.
$ws.Columns.Item(1).ColumnWidth = 8.7
$ws.Columns.Item(2).ColumnWidth = 10.3
$ws.Columns.Item(3).ColumnWidth = 13.6
$ws.Columns.Item(4).ColumnWidth = 12.1
$ws.Columns.Item(5).ColumnWidth = 18.7
$ws.Columns.Item(6).ColumnWidth = 12
$ws.Columns.Item(14).ColumnWidth = 69.8
$ws.Columns.Item(15).ColumnWidth = 7.3
$ws.Columns.Item(16).ColumnWidth = 51.3
$ws.Columns.Item(17).ColumnWidth = 8.6
$ws.Columns.Item(18).ColumnWidth = 4.3
$ws.Columns.Item(19).ColumnWidth = 20.8
